# edit.ps1 - apply "Work flow document" diff via Word COM-interop (InsertXML surgery)
$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Set-ParaXml($paraIndex, $innerXml) {
    $p = $d.Paragraphs($paraIndex)
    $r = $p.Range
    $xml = "<w:p $wns>$innerXml</w:p>"
    $r.InsertXML($xml)
}

# ---------------------------------------------------------------------------
# 1) Title paragraph: drop proofErr gramStart/gramEnd, merge the two runs into one
# ---------------------------------------------------------------------------
$titleInner = '<w:pPr><w:jc w:val="center"/><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:u w:val="single"/></w:rPr><w:t>Work flow document</w:t></w:r>'
Set-ParaXml 1 $titleInner

# ---------------------------------------------------------------------------
# 2) Step 13: merge "Step 13 -Set up a " + "github" + " repository" (drop proofErr spellStart/spellEnd)
# ---------------------------------------------------------------------------
$step13Inner = '<w:r><w:t>Step 13 -Set up a github repository</w:t></w:r><w:r><w:t>.</w:t></w:r>'
Set-ParaXml 15 $step13Inner

# ---------------------------------------------------------------------------
# 3) Step 14: merge "Upload to " + "github" + "." (drop proofErr spellStart/spellEnd)
# ---------------------------------------------------------------------------
$step14Inner = '<w:r><w:t>Step 14</w:t></w:r><w:r><w:t xml:space="preserve"> -</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>Upload to github.</w:t></w:r>'
Set-ParaXml 16 $step14Inner

# ---------------------------------------------------------------------------
# 4) Step 15: merge "Step 15 – download " + "springboot" + "."
# ---------------------------------------------------------------------------
$step15Inner = '<w:r><w:t>Step 15 – download springboot.</w:t></w:r>'
Set-ParaXml 17 $step15Inner

# ---------------------------------------------------------------------------
# 5) Step 16: merge everything into one run
# ---------------------------------------------------------------------------
$step16Inner = '<w:r><w:t>Step 16 – realise springboot is part of intelliJ and unistall springboot.</w:t></w:r>'
Set-ParaXml 18 $step16Inner

# ---------------------------------------------------------------------------
# 6) Step 18 paragraph is fully rewritten
# ---------------------------------------------------------------------------
$step18Inner = '<w:r><w:t xml:space="preserve">Step 18 – </w:t></w:r><w:r><w:t>S</w:t></w:r><w:r><w:t xml:space="preserve">print </w:t></w:r><w:r><w:t>– To create and develop the front end of my project through html/css/javascript</w:t></w:r><w:r><w:t>.</w:t></w:r>'
Set-ParaXml 20 $step18Inner

# ---------------------------------------------------------------------------
# 7) Step 19 paragraph is fully rewritten (with highlight)
# ---------------------------------------------------------------------------
$step19Inner = '<w:r><w:t xml:space="preserve">Step 19 – </w:t></w:r><w:r><w:rPr><w:highlight w:val="darkCyan"/></w:rPr><w:t>S</w:t></w:r><w:r><w:rPr><w:highlight w:val="darkCyan"/></w:rPr><w:t>tand up</w:t></w:r><w:r><w:rPr><w:highlight w:val="darkCyan"/></w:rPr><w:t xml:space="preserve"> – Code cleanliness and structure poor and not enough documentation.</w:t></w:r>'
Set-ParaXml 21 $step19Inner

# ---------------------------------------------------------------------------
# 8) Step 20 paragraph is fully rewritten (bookmark removed from here)
# ---------------------------------------------------------------------------
$step20Inner = '<w:r><w:t>Step 20</w:t></w:r><w:r><w:t>,21,22</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>–</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>Continue to develop front end</w:t></w:r>'
Set-ParaXml 22 $step20Inner

# ---------------------------------------------------------------------------
# 9) Append the new paragraphs (Step 23 .. Step 34), bookmark moves to the very end
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$endPoint = $d.Range($lastPara.Range.End, $lastPara.Range.End)

$newParas = ''
$newParas += "<w:p $wns><w:r><w:t>Step 23 – Compare front end with wire frame and flow chart.</w:t></w:r></w:p>"
$newParas += "<w:p $wns><w:r><w:t>Step 24 – Set up account on google cloud and to upload website onto the cloud platform</w:t></w:r></w:p>"
$newParas += "<w:p $wns><w:r><w:t>Step 2</w:t></w:r><w:r><w:t>5</w:t></w:r><w:r><w:t xml:space=`"preserve`"> – Start to develop back end to meet function in the front end.</w:t></w:r></w:p>"
$newParas += "<w:p $wns><w:r><w:t>Step 2</w:t></w:r><w:r><w:t>6</w:t></w:r><w:r><w:t xml:space=`"preserve`"> – Sprint number 2 </w:t></w:r><w:r><w:t>- A</w:t></w:r><w:r><w:t xml:space=`"preserve`"> save function with the front end, back end and database</w:t></w:r><w:r><w:t xml:space=`"preserve`">. A delete button function that linked with front, back and database. </w:t></w:r><w:r><w:t>Create a table on the page that displayed bookings</w:t></w:r></w:p>"
$newParas += "<w:p $wns><w:r><w:t>Step 2</w:t></w:r><w:r><w:t>7</w:t></w:r><w:r><w:t xml:space=`"preserve`"> – </w:t></w:r><w:r><w:rPr><w:highlight w:val=`"darkCyan`"/></w:rPr><w:t>Stand up – Created the functions that worked however the table wouldn&#8217;t save the input and would vanish off the screen if you refreshed or went onto a different page and then back. Data would however save to the MySql database.</w:t></w:r></w:p>"
$newParas += "<w:p $wns><w:r><w:t>Step 2</w:t></w:r><w:r><w:t>8</w:t></w:r><w:r><w:t xml:space=`"preserve`"> – Review wireframe to check design plans</w:t></w:r></w:p>"
$newParas += "<w:p $wns><w:r><w:lastRenderedPageBreak/><w:t>Step 2</w:t></w:r><w:r><w:t>9</w:t></w:r><w:r><w:t xml:space=`"preserve`"> – To review Kanban board and make adjusts/finish off stories/add new stories.</w:t></w:r></w:p>"
$newParas += "<w:p $wns><w:r><w:t xml:space=`"preserve`">Step </w:t></w:r><w:r><w:t>30</w:t></w:r><w:r><w:t xml:space=`"preserve`"> – To start implementing Tests and test coverage into my program.</w:t></w:r></w:p>"
$newParas += "<w:p $wns><w:r><w:t>Step 31 – Double check progress is uploaded to github and to continue to use Jenkins for CI.</w:t></w:r></w:p>"
$newParas += "<w:p $wns><w:r><w:t>Step 32 – Tidy up code</w:t></w:r></w:p>"
$newParas += "<w:p $wns><w:r><w:t>Step 33 – Start to review documents.</w:t></w:r></w:p>"
$newParas += "<w:p $wns><w:r><w:t>Step 34 – Prep for presentation to nationwide.</w:t></w:r><w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/><w:bookmarkEnd w:id=`"0`"/></w:p>"

$endPoint.InsertXML($newParas)

# ---------------------------------------------------------------------------
# 10) Remove the now-orphaned bookmark left on the old Step 20 paragraph
#     (InsertXML on paragraph 22 above already dropped it, nothing further to do)
# ---------------------------------------------------------------------------
